# Updates the crypto price/volume table (columns D = Price, E = Volume(1h))
# with refreshed figures, and swaps the PEPE/NEARProtocol (rows 25-26) and
# RenderToken/Bittensor (rows 37-38) rank positions.
#
# Price cells (column D) are prefixed with a leading apostrophe so Excel
# stores them as text (matching the source data, which includes
# thousand-separated values like "96.799.43" that are not valid numbers);
# the Style is then reset to "Normal" so the auto-applied Text number
# format doesn't leave a stray style diff behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''96.721.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.79%  '

$ws.Range("E3").Value = '  -0.74%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = '''241.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.50%  '

$ws.Range("D6").Value = '''653.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.70%  '

$ws.Range("E7").Value = '  +14.80%  '

$ws.Range("D8").Value = '''0.413'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.57%  '

$ws.Range("D9").Value = '''1.07'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +8.41%  '

$ws.Range("E10").Value = '  +0.07%  '

$ws.Range("D11").Value = '''3.578.99'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.84%  '

$ws.Range("D12").Value = '''43.51'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.41%  '

$ws.Range("E13").Value = '  +0.73%  '

$ws.Range("E14").Value = '  +1.47%  '

$ws.Range("D15").Value = '''4.248.65'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.05%  '

$ws.Range("D16").Value = '''96.575.91'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.71%  '

$ws.Range("E17").Value = '  +2.81%  '

$ws.Range("D18").Value = '''3.579.22'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.39%  '

$ws.Range("D19").Value = '''8.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.05%  '

$ws.Range("E20").Value = '  -1.24%  '

$ws.Range("D21").Value = '''18.06'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.97%  '

$ws.Range("D22").Value = '''0.539'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +12.54%  '

$ws.Range("D23").Value = '''509.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.29%  '

$ws.Range("E24").Value = '  -3.57%  '

$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").Value = '''0.0000201'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.84%  '

$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").Value = '''6.94'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.58%  '

$ws.Range("D27").Value = '''96.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.22%  '

$ws.Range("E28").Value = '  +5.46%  '

$ws.Range("D29").Value = '''3.774.59'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.67%  '

$ws.Range("D30").Value = '''0.153'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.92%  '

$ws.Range("E31").Value = '  -2.42%  '

$ws.Range("D32").Value = '''11.53'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.97%  '

$ws.Range("E33").Value = '  +0.15%  '

$ws.Range("D34").Value = '''0.185'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.65%  '

$ws.Range("D35").Value = '''0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.74%  '

$ws.Range("D36").Value = '''31.47'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.94%  '

$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").Value = '''624.69'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +10.36%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").Value = '''8.84'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.29%  '

$ws.Range("E39").Value = '  +2.64%  '

$ws.Range("E40").Value = '  +12.26%  '

$ws.Range("E41").Value = '  +1.62%  '

$ws.Range("E42").Value = '  +0.06%  '

$ws.Range("D43").Value = '''0.912'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.37%  '

$ws.Range("E44").Value = '  +5.58%  '

$ws.Range("D45").Value = '''5.83'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.50%  '

$ws.Range("E46").Value = '  +5.04%  '

$ws.Range("E47").Value = '  +3.66%  '

$ws.Range("E48").Value = '  -0.60%  '

$ws.Range("D49").Value = '''32.97'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.54%  '

$ws.Range("D50").Value = '''3.54'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.13%  '

$ws.Range("D51").Value = '''8.31'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.31%  '
